$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append an "order by ... ASC LIMIT 100" clause (matching each tab's natural
# sort key) to the end of the Neo4j query text stored in column B for the
# CasesTab (row 2), SamplesTab (row 3) and FilesTab (row 4) rows.

$filesQuery = $ws.Range("B4").Value2
$ws.Range("B4").Value = $filesQuery + "`n order By f.file_name ASC LIMIT 100"

$samplesQuery = $ws.Range("B3").Value2
$ws.Range("B3").Value = $samplesQuery + "`n order By samp.sample_id ASC LIMIT 100"

$casesQuery = $ws.Range("B2").Value2
$ws.Range("B2").Value = $casesQuery + "`n order By ss.study_subject_id ASC LIMIT 100"

# Match the author's final selection (cell B2 on the CasesTab row).
$ws.Range("B2").Select() | Out-Null

# The extra line in each query re-wraps the cell text, so the rows grow a
# touch taller (values per the saved workbook's recalculated autofit).
$ws.Rows(2).RowHeight = 374.4
$ws.Rows(3).RowHeight = 409.6
$ws.Rows(4).RowHeight = 331.2
